# Apply cell value updates for cryptos.xlsx per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.517.51'
$ws.Range('E2').Value = '  -2.30%  '
$ws.Range('D3').Value = '2.488.29'
$ws.Range('E3').Value = '  -1.26%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '''313.73'
$ws.Range('E5').Value = '  +0.87%  '
$ws.Range('D6').Value = '''94.98'
$ws.Range('E6').Value = '  -3.97%  '
$ws.Range('D7').Value = '''0.549'
$ws.Range('E7').Value = '  -2.57%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -3.15%  '
$ws.Range('E10').Value = '  -4.19%  '
$ws.Range('D11').Value = '''0.0783'
$ws.Range('E11').Value = '  -2.32%  '
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('E13').Value = '  -2.53%  '
$ws.Range('D14').Value = '2.870.22'
$ws.Range('E14').Value = '  -1.27%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.517.81'
$ws.Range('E15').Value = '  +0.14%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = '''15.46'
$ws.Range('E16').Value = '  +1.00%  '
$ws.Range('D17').Value = '''0.796'
$ws.Range('E17').Value = '  -1.17%  '
$ws.Range('D18').Value = '41.488.22'
$ws.Range('E18').Value = '  -2.37%  '
$ws.Range('D19').Value = '''6.36'
$ws.Range('E19').Value = '  -3.57%  '
$ws.Range('D20').Value = '0.0₃0926'
$ws.Range('E20').Value = '  -1.92%  '
$ws.Range('D21').Value = '''11.35'
$ws.Range('E21').Value = '  -6.07%  '
$ws.Range('D22').Value = '''69.04'
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('D23').Value = '''237.63'
$ws.Range('E23').Value = '  -1.29%  '
$ws.Range('D24').Value = '''2.77'
$ws.Range('E24').Value = '  -2.40%  '
$ws.Range('D25').Value = '''1.91'
$ws.Range('E25').Value = '  -3.92%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').Value = '''24.27'
$ws.Range('E27').Value = '  -4.15%  '
$ws.Range('D28').Value = '''2.24'
$ws.Range('E28').Value = '  -1.07%  '
$ws.Range('D29').Value = '''9.81'
$ws.Range('E29').Value = '  -2.14%  '
$ws.Range('D30').Value = '''36.83'
$ws.Range('E30').Value = '  -3.11%  '
$ws.Range('D31').Value = '''152.73'
$ws.Range('E31').Value = '  -2.49%  '
$ws.Range('D32').Value = '''5.52'
$ws.Range('E32').Value = '  -5.03%  '
$ws.Range('E33').Value = '  -3.57%  '
$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D34').Value = '''18.19'
$ws.Range('E34').Value = '  +4.08%  '
$ws.Range('D35').Value = '''0.0758'
$ws.Range('E35').Value = '  -3.47%  '
$ws.Range('B36').Value = 'ApeXProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D36').Value = '''2.52'
$ws.Range('E36').Value = '  -6.84%  '
$ws.Range('D37').Value = '''3.10'
$ws.Range('E37').Value = '  -1.51%  '
$ws.Range('D38').Value = '''1.89'
$ws.Range('E38').Value = '  -3.16%  '
$ws.Range('E39').Value = '  -1.61%  '
$ws.Range('E40').Value = '  -6.95%  '
$ws.Range('D41').Value = '''4.22'
$ws.Range('E41').Value = '  +2.06%  '
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('D43').Value = '''19.75'
$ws.Range('E43').Value = '  -9.91%  '
$ws.Range('D44').Value = '2.000.36'
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('E45').Value = '  -2.87%  '
$ws.Range('E46').Value = '  -6.58%  '
$ws.Range('D47').Value = '''8.85'
$ws.Range('E47').Value = '  -1.42%  '
$ws.Range('D48').Value = '2.733.98'
$ws.Range('E48').Value = '  -0.87%  '
$ws.Range('D49').Value = '''70.04'
$ws.Range('E49').Value = '  -1.98%  '
$ws.Range('D50').Value = '''97.42'
$ws.Range('E50').Value = '  -3.05%  '
$ws.Range('D51').Value = '''0.179'
$ws.Range('E51').Value = '  -5.43%  '
